$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 157.6633739787549
$ws.Range("D2").Value = 782.0969401757909
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = 203.0269105549284
$ws.Range("C3").Value = 8057
